$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Build the "Times New Roman 12" data font on a scratch cell, then copy
#    just the format onto the whole existing table body (A1:C45). Doing the
#    property writes on a single scratch cell (rather than the 135-cell
#    range directly) keeps the resulting style table clean (one new cellXf
#    instead of one per touched cell).
# ---------------------------------------------------------------------------
$ws.Range("Z1").Font.Size = 12
$ws.Range("Z1").Font.Name = "Times New Roman"
$ws.Range("Z1").Copy()
$ws.Range("A1:C45").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Build the big bold title font + centered alignment on another scratch
#    cell.
# ---------------------------------------------------------------------------
$ws.Range("Z2").Font.Size = 36
$ws.Range("Z2").Font.Name = "Times New Roman"
$ws.Range("Z2").Font.Bold = $true
$ws.Range("Z2").HorizontalAlignment = -4108
$ws.Range("Z2").VerticalAlignment = -4108

# Clean up the scratch cells now that their formats have been harvested.
$ws.Range("Z1:Z2").Clear()

# ---------------------------------------------------------------------------
# 3) Insert 3 rows above the current header row, pushing the header from
#    row 1 -> row 4 and the data from rows 2:45 -> rows 5:48.
# ---------------------------------------------------------------------------
$ws.Range("A1:C3").EntireRow.Insert()

# 4) Resize the table (ListObject) to the new location.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A4:C48"))

# ---------------------------------------------------------------------------
# 5) Write & format the title in A2, merged across A2:C2, re-using the title
#    font/alignment we built above (re-create it directly since the scratch
#    cell is gone, then copy it down).
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "Dictionnaire des données"
$ws.Range("A2:C2").Merge()

$ws.Range("Z1").Font.Size = 36
$ws.Range("Z1").Font.Name = "Times New Roman"
$ws.Range("Z1").Font.Bold = $true
$ws.Range("Z1").HorizontalAlignment = -4108
$ws.Range("Z1").VerticalAlignment = -4108
$ws.Range("Z1").Copy()
$ws.Range("A2:C2").PasteSpecial(-4122)
$ws.Range("Z1").Clear()

$ws.Range("A2:C2").RowHeight = 45

# 6) Column widths.
$ws.Columns("A").ColumnWidth = 19.28515625
$ws.Columns("B").ColumnWidth = 100.42578125
$ws.Columns("C").ColumnWidth = 13.85546875

# 7) Page setup: fit to page, paper size A4, 65% scale, portrait.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
$ws.PageSetup.Zoom = 65
$ws.PageSetup.FitToPagesWide = $false
$ws.PageSetup.FitToPagesTall = $false

# 8) Selection/scroll matching the author's final view.
$ws.Range("A2:C48").Select()
